$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 256
$ws.Range("B2").Value = 238385
$ws.Range("C2").Value = 41727
$ws.Range("D2").Value = 341388
$ws.Range("E2").Value = 381541
$ws.Range("F2").Value = 426050
$ws.Range("G2").Value = 406224
$ws.Range("H2").Value = 501930
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()

# Row 5 updates
$ws.Range("A5").Value = 256
$ws.Range("B5").Value = 390846
$ws.Range("C5").Value = 665175
$ws.Range("D5").Value = 63008
$ws.Range("E5").Value = 920700
$ws.Range("F5").Value = 853379
$ws.Range("G5").Value = 810234
$ws.Range("H5").Value = 815155
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
